# Generate Report for Handback
#
# Applies the "handback" report-generation update to localization-status.xlsx:
#   - Overview: status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - zh-cn / de-de detail sheets: the "Latest Target File" / "Latest Handback
#     File" / "Latest Handback DateTime" columns get populated now that the
#     handback round-trip finished (they were blank / epoch placeholders).
#   - Column widths are widened on the columns that now hold long file names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: the engine quantizes ColumnWidth to 1/6-character increments when
# converting to the stored <col width> attribute (stored = ColumnWidth +
# 5/6, rounded to the nearest 1/6). Back the desired stored width out to the
# ColumnWidth value that reproduces it as closely as possible.
# ---------------------------------------------------------------------------
function Set-StoredColumnWidth($range, [double]$storedWidth) {
    $range.ColumnWidth = $storedWidth - 0.8333333333333334
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

Set-StoredColumnWidth $wsOverview.Range("E:F") 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2 -> 409ba8ab-...md handback
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51610b302e7c3fd99f3fe972a096e4c2a6906b57/e2e/409ba8ab-6a1c-4494-9e24-d3cd1508d24c.md",
    "",
    "",
    "409ba8ab-6a1c-4494-9e24-d3cd1508d24c.md"
) | Out-Null
$wsZhCn.Range("J2").Value = "409ba8ab-6a1c-4494-9e24-d3cd1508d24c.3929a4c46c2ddf1cc18aa338ed09182faf64de95.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-04 10:52:37"

# Row 3 -> 5702a499-...md handback
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51610b302e7c3fd99f3fe972a096e4c2a6906b57/e2e/5702a499-1e12-4816-b4f0-8e0f0822bf6e.md",
    "",
    "",
    "5702a499-1e12-4816-b4f0-8e0f0822bf6e.md"
) | Out-Null
$wsZhCn.Range("J3").Value = "5702a499-1e12-4816-b4f0-8e0f0822bf6e.5ee0210c63bc3fd7e43adf6de0b13d2c142ad1c5.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-04 10:52:37"

Set-StoredColumnWidth $wsZhCn.Range("C:C") 29.9777047293527
Set-StoredColumnWidth $wsZhCn.Range("I:J") 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2 -> 409ba8ab-...md handback
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51610b302e7c3fd99f3fe972a096e4c2a6906b57/e2e/409ba8ab-6a1c-4494-9e24-d3cd1508d24c.md",
    "",
    "",
    "409ba8ab-6a1c-4494-9e24-d3cd1508d24c.md"
) | Out-Null
$wsDeDe.Range("J2").Value = "409ba8ab-6a1c-4494-9e24-d3cd1508d24c.3929a4c46c2ddf1cc18aa338ed09182faf64de95.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-04 10:52:44"

# Row 3 -> 5702a499-...md handback
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51610b302e7c3fd99f3fe972a096e4c2a6906b57/e2e/5702a499-1e12-4816-b4f0-8e0f0822bf6e.md",
    "",
    "",
    "5702a499-1e12-4816-b4f0-8e0f0822bf6e.md"
) | Out-Null
$wsDeDe.Range("J3").Value = "5702a499-1e12-4816-b4f0-8e0f0822bf6e.5ee0210c63bc3fd7e43adf6de0b13d2c142ad1c5.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-04 10:52:44"

Set-StoredColumnWidth $wsDeDe.Range("C:C") 29.9777047293527
Set-StoredColumnWidth $wsDeDe.Range("I:J") 40

Write-Host "Handback report generated."
